$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate column P (years 2021 data) into new column Q as the starting
# point, carrying over both the cell styles and values.
$ws.Range("P3:P25").Copy($ws.Range("Q3:Q25"))

# Update the new column's header (year) and the data that changed for 2022.
$ws.Range("Q4").Value = 2022
$ws.Range("Q5").Value = 8725
$ws.Range("Q7").Value = 8347
$ws.Range("Q8").Value = 378

# Move the active selection to Q3, matching the refreshed sheet view.
$ws.Range("Q3").Select()
